# ---------------------------------------------------------------------------
# Edit: mlk.docx
#   1. Append two trailing spaces to the first paragraph's existing text,
#      then append a parenthetical note in red (C00000), split across three
#      runs exactly as captured by the original authoring session.
#   2. Append a new, otherwise-empty paragraph at the very end of the body
#      (before the final section break) shaded with fill F9F9F9.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. First paragraph: "This is a Microsoft word document." -------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1                 # exclude the paragraph mark
$r1.InsertAfter("  ")                 # two trailing spaces, same formatting

$p1 = $d.Paragraphs(1)
$r2 = $p1.Range
$r2.End = $r2.End - 1
$r2.Collapse(0)
$r2.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r2.Font.Color = 192                  # wdColor BGR for RGB C00000

$p1 = $d.Paragraphs(1)
$r3 = $p1.Range
$r3.End = $r3.End - 1
$r3.Collapse(0)
$r3.InsertAfter("rsion for branch alternate")
$r3.Font.Color = 192

$p1 = $d.Paragraphs(1)
$r4 = $p1.Range
$r4.End = $r4.End - 1
$r4.Collapse(0)
$r4.InsertAfter(")")
$r4.Font.Color = 192

# --- 2. New shaded paragraph at the end of the document --------------------
$tail = $d.Content
$tail.Collapse(0)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$tail.InsertXML($newParaXml)
